$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value, reflecting the refreshed coinranking.com
# export (prices, 1h volume deltas, and the Cosmos/EthereumClassic and
# WEMIXToken/BabyDogeCoin rows swapping rank position).
$updates = [ordered]@{
    "D2"  = "33.855.83"
    "E2"  = "  -1.05%  "
    "D3"  = "1.779.63"
    "E3"  = "  -1.42%  "
    "E4"  = "  +0.14%  "
    "D5"  = "224.27"
    "E5"  = "  +0.30%  "
    "E6"  = "  -1.47%  "
    "E7"  = "  +0.09%  "
    "D8"  = "31.61"
    "E8"  = "  -4.46%  "
    "D9"  = "0.289"
    "E9"  = "  +0.15%  "
    "D10" = "0.0677"
    "E10" = "  -5.99%  "
    "D11" = "0.0935"
    "E11" = "  +0.67%  "
    "D12" = "2.029.62"
    "E12" = "  -1.68%  "
    "E13" = "  +0.55%  "
    "D14" = "1.788.02"
    "E14" = "  -1.09%  "
    "D15" = "33.859.47"
    "E15" = "  -1.03%  "
    "D16" = "0.609"
    "E16" = "  -3.84%  "
    "E17" = "  -2.56%  "
    "D18" = "66.76"
    "E18" = "  -2.96%  "
    "D19" = "238.84"
    "E20" = "  -2.37%  "
    "E21" = "  +0.19%  "
    "D22" = "10.54"
    "E22" = "  -4.76%  "
    "E23" = "  -2.76%  "
    "E24" = "  -3.25%  "
    "D25" = "161.03"
    "E25" = "  +0.84%  "
    "B26" = "EthereumClassic"
    "C26" = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
    "D26" = "16.08"
    "E26" = "  -3.44%  "
    "B27" = "Cosmos"
    "C27" = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
    "D27" = "7.02"
    "E27" = "  -2.05%  "
    "E28" = "  -1.24%  "
    "E29" = "  +0.23%  "
    "E30" = "  +0.61%  "
    "D31" = "0.0509"
    "E31" = "  -3.91%  "
    "E32" = "  -4.37%  "
    "E33" = "  -0.54%  "
    "E34" = "  -2.54%  "
    "D35" = "1.387.83"
    "E35" = "  -2.23%  "
    "D36" = "0.636"
    "E36" = "  -2.61%  "
    "E37" = "  -1.76%  "
    "E38" = "  -1.39%  "
    "D39" = "2.26"
    "E39" = "  +4.02%  "
    "E40" = "  -0.01%  "
    "D41" = "0.913"
    "E41" = "  -3.64%  "
    "D42" = "78.28"
    "E42" = "  -3.08%  "
    "E43" = "  -3.12%  "
    "D44" = "13.48"
    "E44" = "  +11.50%  "
    "D45" = "0.0511"
    "E45" = "  +2.73%  "
    "B46" = "WEMIXToken"
    "C46" = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
    "D46" = "1.07"
    "E46" = "  +2.26%  "
    "B47" = "BabyDogeCoin"
    "C47" = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
    "D47" = "0.0₆0135"
    "E47" = "  +8.45%  "
    "E48" = "  -1.91%  "
    "D49" = "1.935.83"
    "E49" = "  -1.49%  "
    "D50" = "105.19"
    "E50" = "  -3.00%  "
    "E51" = "  +0.08%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $value = $updates[$ref]

    # Price cells in column D such as "224.27" or "31.61" look like plain
    # decimal numbers to Excel. Assigning them straight to .Value would let
    # Excel silently reinterpret them as floating point numbers (losing the
    # exact source formatting, e.g. "224.27" -> 224.27000000000001). Detect
    # that case and force the cell to Text first so the literal string is
    # kept, then drop the style back to Normal so no stray number format is
    # left behind on the cell.
    if ($ref.Substring(0,1) -eq "D" -and $value -match '^[0-9]+\.[0-9]+$') {
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}
